$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Add time_taken values for rows 2-6
$ws.Range("F2").Value = "2021-10-05 10:52:31.999192"
$ws.Range("F3").Value = "2021-10-05 10:52:31.999204"
$ws.Range("F4").Value = "2021-10-05 10:52:31.999208"
$ws.Range("F5").Value = "2021-10-05 10:52:31.999211"
$ws.Range("F6").Value = "2021-10-05 10:52:31.999215"
